# Update the "取得日時" (acquired datetime) column on the ランサーズ sheet
# for rows 2-7 to reflect the latest fetch timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-20 06:30:31"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
